# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values replacing the previous "Strike#" counts in column G
$kValues = @{
    2  = 8
    3  = 9
    4  = 5
    5  = 6
    6  = 12
    7  = 4
    8  = 7
    9  = 5
    10 = 8
    11 = 11
    12 = 5
    13 = 8
    14 = 6
    15 = 13
    16 = 5
    17 = 7
    18 = 3
    19 = 7
    20 = 5
    21 = 6
    22 = 4
    23 = 0
    24 = 6
    25 = 9
    26 = 3
    27 = 7
    28 = 3
    29 = 4
    30 = 9
    31 = 9
    32 = 7
    33 = 4
    34 = 2
    35 = 5
    36 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
